$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("latest")

$values = @{
    2  = @(0.186985831112528, 0.722203634391301, 0.9341118818059831, 0.9664946362013517, 0.9757249869846261, 18)
    3  = @(0.2330935088154202, 0.4998167590291169, 0.3374744503351871, 0.580925511864634, 0.5484872224022262, 17)
    4  = @(0.2474022870896494, 0.4562231751086744, 0.282823943946844, 0.531811944155868, 0.4862000162233972, 16)
    5  = @(0.3331709995737225, 0.5061089045982429, 0.3514486812937936, 0.5928310731513604, 0.5075632610151456, 15)
    6  = @(0.3697548977873877, 0.5076571467308241, 0.3515889833503586, 0.5929493935829251, 0.4810392743262605, 14)
    7  = @(0.3550099532386083, 0.5235044444458751, 0.3739286097767974, 0.611497023522435, 0.518222527605833, 13)
    8  = @(0.4180599993501528, 0.5377036271313899, 0.3917271539151825, 0.6258811020594747, 0.4864935662755369, 12)
    9  = @(0.4373552979544756, 0.5644811036562025, 0.414358809081934, 0.643707083293274, 0.4953655899196011, 11)
    10 = @(0.4113383373941182, 0.553516800113468, 0.4096270470205588, 0.640021130135997, 0.5168578347523529, 10)
    11 = @(0.3463436447865475, 0.5120086005850197, 0.3587015581898463, 0.5989169877285552, 0.5182577473073058, 9)
}

$cols = @("B", "C", "D", "E", "F", "G")

foreach ($row in $values.Keys) {
    $rowVals = $values[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $rowVals[$i]
    }
}
